$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in cell A6: "00000998" -> "00000988"
$ws.Range("A6").Value = "00000988"

# Move the active selection to A7, as if the user clicked below the table
$ws.Range("A7").Select()
